$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for the new "experience" columns ---
# (inputs chosen so the engine's stored `width` lands as close as
# possible to the authored 16.6640625 / 13.9140625 target values)
$ws.Columns.Item(7).ColumnWidth = 15.833333333333334
$ws.Columns.Item(8).ColumnWidth = 13

# --- Write the values in the same order the shared-string table was
#     originally built in (H2, G1, G2, I1, I2, H1, then the rest) so the
#     resulting sharedStrings.xml ordering lines up with the target file.
$ws.Range("H2").Value = "January 1,2024"
$ws.Range("G1").Value = "Start Date"
$ws.Range("G2").Value = "January 1,2022"
$ws.Range("I1").Value = "Position"
$ws.Range("I2").Value = "Software Engineer I"
$ws.Range("H1").Value = "End Date"

$ws.Range("H3").Value = "January 1,2024"
$ws.Range("G3").Value = "January 1,2022"
$ws.Range("I3").Value = "Software Engineer I"

$ws.Range("H4").Value = "January 1,2024"
$ws.Range("G4").Value = "January 1,2022"
$ws.Range("I4").Value = "Software Engineer I"

$ws.Range("H5").Value = "January 1,2024"
$ws.Range("G5").Value = "January 1,2022"
$ws.Range("I5").Value = "Software Engineer I"

$ws.Range("H6").Value = "January 1,2024"
$ws.Range("G6").Value = "January 1,2022"
$ws.Range("I6").Value = "Software Engineer I"

# --- Borders: thin left/right border around the new header cells
#     (applied cell-by-cell so each one keeps both its own left and
#     right edge rather than only the outer edge of a merged
#     selection). ---
foreach ($addr in @("G1", "H1", "I1")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).ColorIndex = -4105
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).ColorIndex = -4105
}

# --- Selection / active cell, matches the state left behind in the
#     authored workbook. ---
$ws.Range("H14").Select()

Write-Output "experience columns added"
